# "Added homework for 05. View-Engines"
# The canonical diff for this commit removes the trailing "Homework" slide
# (slide 33 - Title "Homework" / Content Placeholder with the assignment
# bullet list) from the deck: it drops the <p:sldId id="286" r:id="rId34"/>
# entry from the presentation's slide list and deletes the slide part
# itself. Deleting the slide through the Slides collection takes care of
# both (PowerPoint removes the slide's <p:sldId> entry from sldIdLst and
# drops the slide part/relationship automatically).

$p = $ppt.ActivePresentation

# The Homework slide is the last slide (33) in the deck.
$p.Slides.Item($p.Slides.Count).Delete()
